$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark so it starts at the very beginning of the
#    document (before the "Letter " run) instead of between "Letter " and
#    "of Invitation". Re-adding a bookmark with the same name relocates it
#    (removing the old bookmarkStart/bookmarkEnd pair and inserting a new one
#    spanning the given range), and Word always collapses the bookmarkEnd to
#    sit right after the covered text - exactly where the old bookmark used
#    to live.
# ---------------------------------------------------------------------------
$goBackRange = $d.Range(0, 7)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $old with $new, then "touch" the
# formatting of the freshly written text (bold on/off, a no-op visually) so
# the engine keeps it as its own run instead of silently re-merging it with
# an identically formatted neighbouring run.
# ---------------------------------------------------------------------------
function Replace-KeepOwnRun([string]$old, [string]$new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return $false }
    $start = $rng.Start
    $rng.Text = $new
    $newRange = $d.Range($start, $start + $new.Length)
    $newRange.Font.Bold = $true
    $newRange.Font.Bold = $false
    return $true
}

# ---------------------------------------------------------------------------
# 2) "Middle Name" -> "Middle Name(s)" (keep "Middle" / " Name" runs intact,
#    add a brand-new "(s)" run right before the superscript "1" footnote
#    marker).
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Middle Name") | Out-Null
$findRange.Collapse(0)  # wdCollapseEnd
$findRange.InsertBefore("(s)")
$findRange.Font.Bold = $true
$findRange.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) "Cell Telephone Number" -> "Mobile Number"
#    (merge the "Cell" / " " / "Telephone" runs into a single "Mobile" run,
#    leaving the trailing " Number" run untouched).
# ---------------------------------------------------------------------------
Replace-KeepOwnRun "Cell Telephone" "Mobile" | Out-Null

# ---------------------------------------------------------------------------
# 4) "Surname and First Name must match..." -> "Complete name must match..."
#    (merge the "Surname and F" / "irst " / "Name" runs into one "Complete
#    name" run, leaving the trailing " must match the name on your
#    passport." run untouched).
# ---------------------------------------------------------------------------
Replace-KeepOwnRun "Surname and First Name" "Complete name" | Out-Null
